$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 343, pushing the existing
# rows 343-350 down to 345-352.
$ws.Rows("343:344").Insert()

# New row 343: Papa / Asterix, 2022-02-03 (44595), Region del Maule
$ws.Cells.Item(343, 1).Value = 8
$ws.Cells.Item(343, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(343, 3).Value = "Coquimbo"
$ws.Cells.Item(343, 4).Value = 44595
$ws.Cells.Item(343, 5).Value = 4
$ws.Cells.Item(343, 6).Value = 100114001
$ws.Cells.Item(343, 7).Value = "Papa"
$ws.Cells.Item(343, 8).Value = "Asterix"
$ws.Cells.Item(343, 9).Value = "1a (cosecha)"
$ws.Cells.Item(343, 10).Value = 2400
$ws.Cells.Item(343, 11).Value = 9500
$ws.Cells.Item(343, 12).Value = 10000
$ws.Cells.Item(343, 13).Value = 9750
$ws.Cells.Item(343, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(343, 15).Value = "Región del Maule"
$ws.Cells.Item(343, 16).Value = 390
$ws.Cells.Item(343, 17).Value = 25
$ws.Cells.Item(343, 18).Value = "Hortaliza"

# New row 344: Papa / Cardinal, 2022-02-03 (44595), Provincia del Elqui
$ws.Cells.Item(344, 1).Value = 8
$ws.Cells.Item(344, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(344, 3).Value = "Coquimbo"
$ws.Cells.Item(344, 4).Value = 44595
$ws.Cells.Item(344, 5).Value = 4
$ws.Cells.Item(344, 6).Value = 100114001
$ws.Cells.Item(344, 7).Value = "Papa"
$ws.Cells.Item(344, 8).Value = "Cardinal"
$ws.Cells.Item(344, 9).Value = "1a (cosecha)"
$ws.Cells.Item(344, 10).Value = 2000
$ws.Cells.Item(344, 11).Value = 11000
$ws.Cells.Item(344, 12).Value = 12000
$ws.Cells.Item(344, 13).Value = 11500
$ws.Cells.Item(344, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(344, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(344, 16).Value = 460
$ws.Cells.Item(344, 17).Value = 25
$ws.Cells.Item(344, 18).Value = "Hortaliza"
